# Auto-generated edit script applying the diff to 广州-漫展信息.xlsx
# Updates 'F' column (想去人数 / interested-count) numbers and two stale
# cover-image URLs ('I4' on 演出, 'I13' on 全部类型) per the commit
# "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 335
$ws.Range("F4").Value = 1317
$ws.Range("F5").Value = 378
$ws.Range("F6").Value = 358
$ws.Range("F7").Value = 3914
$ws.Range("F8").Value = 240
$ws.Range("F9").Value = 778
$ws.Range("F10").Value = 2334
$ws.Range("F11").Value = 354
$ws.Range("F13").Value = 754
$ws.Range("F14").Value = 199
$ws.Range("F15").Value = 188
$ws.Range("F16").Value = 2254
$ws.Range("F18").Value = 29
$ws.Range("F21").Value = 236

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 52
$ws.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202410/mx4x20ma1729586170002.jpeg"
$ws.Range("F12").Value = 7
$ws.Range("F22").Value = 65

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6401
$ws.Range("F4").Value = 2116
$ws.Range("F5").Value = 344
$ws.Range("F6").Value = 13

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6401
$ws.Range("F4").Value = 2116
$ws.Range("F5").Value = 344
$ws.Range("F6").Value = 52
$ws.Range("F7").Value = 52
$ws.Range("F8").Value = 5
$ws.Range("F10").Value = 335
$ws.Range("F11").Value = 1317
$ws.Range("F12").Value = 378
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202410/mx4x20ma1729586170002.jpeg"
$ws.Range("F16").Value = 13
$ws.Range("F17").Value = 358
$ws.Range("F18").Value = 3914
$ws.Range("F20").Value = 240
$ws.Range("F24").Value = 778
$ws.Range("F25").Value = 2334
$ws.Range("F26").Value = 354
$ws.Range("F29").Value = 754
$ws.Range("F30").Value = 199
$ws.Range("F31").Value = 188
$ws.Range("F32").Value = 7
$ws.Range("F34").Value = 2254
$ws.Range("F38").Value = 29
$ws.Range("F41").Value = 236
$ws.Range("F49").Value = 65

